# Auto update Excel log
# Appends newly-logged sensor readings to the PIR, Humidity and Temperature
# sheets of the SeniorConnect master log.

$wb = $excel.ActiveWorkbook

$theDate = "2026-01-28"

function Add-LogRows {
    param(
        $ws,
        $rows
    )

    # Force column A (Date) and column E (Value) to Text format before
    # writing so Excel does not "helpfully" reinterpret strings like
    # "2026-01-28" as a real date, or "87.2%" as a numeric percentage.
    $firstRow = $rows[0][0]
    $lastRow = $rows[$rows.Length - 1][0]
    $ws.Range("A$firstRow`:A$lastRow").NumberFormat = "@"
    $ws.Range("E$firstRow`:E$lastRow").NumberFormat = "@"

    foreach ($row in $rows) {
        $r = $row[0]
        $ws.Cells.Item($r, 1).Value = $theDate
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).Value = $row[4]
        $ws.Cells.Item($r, 6).Value = $row[5]
    }
}

# ---------------------------------------------------------------------------
# PIR sheet: rows 148-160
# ---------------------------------------------------------------------------
$wsPIR = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @(148, "14:59:55", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(149, "14:59:56", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(150, "15:00:00", "15:00", "Bathroom", "No Motion", "Inactive"),
    @(151, "15:00:04", "15:00", "Bathroom", "No Motion", "Inactive"),
    @(152, "15:00:10", "15:00", "Bathroom", "No Motion", "Inactive"),
    @(153, "15:00:15", "15:00", "Bathroom", "No Motion", "Inactive"),
    @(154, "15:00:20", "15:00", "Bathroom", "No Motion", "Inactive"),
    @(155, "15:00:25", "15:00", "Bathroom", "No Motion", "Inactive"),
    @(156, "15:00:30", "15:00", "Bathroom", "No Motion", "Inactive"),
    @(157, "15:00:35", "15:00", "Bathroom", "No Motion", "Inactive"),
    @(158, "15:00:40", "15:00", "Bathroom", "No Motion", "Inactive"),
    @(159, "15:00:45", "15:00", "Bathroom", "No Motion", "Inactive"),
    @(160, "15:00:50", "15:00", "Bathroom", "No Motion", "Inactive")
)
Add-LogRows $wsPIR $pirRows

# ---------------------------------------------------------------------------
# Humidity sheet: rows 144-154
# ---------------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @(144, "14:59:55", "14:00", "Bathroom", "87.2%", "Active"),
    @(145, "14:59:58", "14:00", "Bathroom", "88.1%", "Active"),
    @(146, "15:00:06", "15:00", "Bathroom", "87.3%", "Active"),
    @(147, "15:00:18", "15:00", "Bathroom", "88.1%", "Active"),
    @(148, "15:00:22", "15:00", "Bathroom", "88.2%", "Active"),
    @(149, "15:00:26", "15:00", "Bathroom", "87.2%", "Active"),
    @(150, "15:00:30", "15:00", "Bathroom", "88.2%", "Active"),
    @(151, "15:00:38", "15:00", "Bathroom", "87.3%", "Active"),
    @(152, "15:00:42", "15:00", "Bathroom", "88.2%", "Active"),
    @(153, "15:00:50", "15:00", "Bathroom", "88.2%", "Active"),
    @(154, "15:00:54", "15:00", "Bathroom", "88.2%", "Active")
)
Add-LogRows $wsHumidity $humidityRows

# ---------------------------------------------------------------------------
# Temperature sheet: rows 144-154
# ---------------------------------------------------------------------------
$wsTemperature = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @(144, "14:59:56", "14:00", "Bathroom", "22.8C", "Active"),
    @(145, "14:59:58", "14:00", "Bathroom", "22.8C", "Active"),
    @(146, "15:00:06", "15:00", "Bathroom", "22.9C", "Active"),
    @(147, "15:00:19", "15:00", "Bathroom", "22.8C", "Active"),
    @(148, "15:00:23", "15:00", "Bathroom", "22.9C", "Active"),
    @(149, "15:00:27", "15:00", "Bathroom", "22.8C", "Active"),
    @(150, "15:00:31", "15:00", "Bathroom", "22.9C", "Active"),
    @(151, "15:00:39", "15:00", "Bathroom", "22.9C", "Active"),
    @(152, "15:00:43", "15:00", "Bathroom", "22.9C", "Active"),
    @(153, "15:00:51", "15:00", "Bathroom", "22.9C", "Active"),
    @(154, "15:00:55", "15:00", "Bathroom", "22.9C", "Active")
)
Add-LogRows $wsTemperature $temperatureRows

$wb.Save()
